# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

$wsMap = @{}
foreach ($s in $wb.Worksheets) {
    $wsMap[$s.Name] = $s
}

# Update per-play logged data strings (YDS and ST sheets)
$wsMap['YDS'].Range('B2').Value = '5 15 -1 5 9 -2 10 4 -1 1 2 1 4 4 10 -1 2 2 1 0 3 5 5 3 5 1 10 3 1 1 3 8 1 3 5 6 -2 7 4 2 1 3 1 6 9 4 -2 5 -2 4 6 4 -2 0 3 2 5 7 1 4 7 35 11 1 6 2 16 10 1 2 4 3 16 3 4 1 1 4 0 4 4 8 2 0 1 -2 -3 3 2 5 5 0 13 5 3 0 4 6 1 3 4 11 8 7 35 4 4 4 7 24 14 1 3 0 5 4 2 7 11 6 3 1 2 4 9 1 -4 6 8 -1 5 4 -2 0 8 4 9 -3 -4 2 -2 2 4 4 7 4 8 -1 2 4 3 6 -2 0 3 9 0 0 4 1 4 1 2 2 0 0 0 1 2 1 2 3 1 2 5 13 2 9 -6 -3 3 8 1 2 10 1 1 4 2 0 12 8 4 3 1 1 13 9 0 2 4 -4 11 2 2 1 3 10 14 4 1 6 8 3 2 1 3 -4 5 0 4 13 4 2 1 3 8 4 -3 4 5 8 -4 -2 3 2 6 12 3 3 3 2 1 4 3 0 4 0 1 2 3 6 7 4 3 4 3 2 3 3 0 5 7 5 3 1 0 5 3 4 -1 3 7 2 3 -4 4 7 0 1 11 2 5 24 3 8 3 4 22 -2 4 2 2 2 6 4 2 15 -1 7 5 1 2 1 0 2 0 7 13 1 8 0 -7 3 2 8 4 0 3 8 1 3 1 3 5 6 5 6 3 3 -1 7 0 -1 5 1 16 -4 3 3 0 3 2 7 -2 4 1 8 14 4 3 4 -7 5 0 6 3 4 6 6 3 3 2 2 1 5 0 5 0 4 0 4 2 1 4 0 1 2 2 -1 -6 9 3 11 4 9 62 1 2 4 12 4 3 -2 3 10 11 14 1 8 7 12 3 2 11 3 1 4 3 -1 0 0 -1 3 8 5 7 -1 0 2 3 10 0 0 8 2 0 8 1 -3 1 8 6 0 3 4 1 0 5 2 11 1 3 11 5 10 3 3 3 2 0 2 1 3 -1 5 2 2 8 2 3 2 4 -3 10 2 5 -3 1 2 -3 3 3 13 5 2 0 13 1 2 -3 4 17 10 1 6 1 3 9 0 3 1 4 0 8 4 3 17 0 5 2 1 4 1 1 3 -2 3 6 3 1 7 3 13 1 2 7 2 2 2 2 1 3 11 2 2 3 5 3 2 2 4 1 -1 9 0 5 2 18 2 8 -2 3 3 14 5 0 5 3 1 1 0 4 1 5 0 2 2 1 3 4 4 5 1 2 1 0 1 3 -4 -1 2 6 14 3 2 3 5 -3 4 1 -5 1 6 2 4 7 17 1 3 -3 3 0 12 2 1 0 2 5 5 0 -1 0 7 5 -2 4 19 15 1 11 7 1 27 4 12 4 2 1 2 7 0 6 1 2 6 2 5 2 7 39 17 6 5 1 5 11 0 4 2 0 4 0 3 1 0 2 3 0 14 0 5 -1 0 2 6 5 4 11 0 4 -1 4 11 4 0 11 1 4 3 11 8 5 9 0 2 1 -2 1 11 11 4 2 3 5 3 4 0 0 7 2 7 1 3 2 3 3 6 0 0 2 -1 5 4 8 2 1 3 3 2 2 3 2 1 0 8 3 2 9 3 3 -2 3 0'
$wsMap['YDS'].Range('B3').Value = '19 12 28 0 27 15 16 11 4 5 13 10 4 44 1 11 7 19 8 12 18 18 15 13 7 4 24 4 10 12 1 2 12 -1 7 16 22 2 42 4 25 5 4 3 14 11 24 1 10 2 6 7 22 2 8 4 12 14 10 63 1 16 2 9 14 6 12 5 2 13 2 24 13 22 3 12 13 6 8 5 13 14 6 1 9 4 10 9 18 16 6 6 11 4 16 4 6 13 11 19 14 21 7 7 8 4 2 13 7 8 -1 20 6 9 6 8 8 17 4 26 13 26 11 9 8 3 20 23 8 5 17 3 3 18 14 9 4 5 8 3 17 4 3 17 7 27 7 6 7 14 31 8 5 12 -2 7 16 7 9 19 11 27 11 17 4 13 13 4 9 5 17 8 8 5 28 7 12 6 9 8 16 16 24 28 11 7 13 20 16 18 8 19 19 7 7 23 6 0 12 9 9 5 20 11 5 5 -4 7 8 6 3 51 7 5 11 42 9 6 12 11 21 6 21 10 24 8 -4 3 3 46 10 4 1 17 8 11 17 7 11 8 4 3 8 18 9 11 16 0 11 12 6 17 12 14 2 8 8 10 14 6 12 11 3 7 4 -3 4 2 12 20 -2 35 3 9 7 31 10 18 24 14 9 10 27 18 22 4 8 4 5 39 35 7 7 32 19 9 -2 22 4 7 12 9 0 14 8 3 17 2 10 4 23 27 5 3 25 5 5 8 13 14 12 2 7 14 18 5 15 12 -1 13 28 4 14 7 4 -4 11 33 -2 26 -1 3 15 7 7 5 8 54 5 7 0 15 17 6 9 6 2 7 6 12 7 26 4 13 13 7 4 8 2 11 10 16 6 7 5 15 6 17 11 19 14 13 0 3 11 18 -1 16 1 11 12 7 3 7 9 9 5 18 5 12 6 9 7 6 8 7 5 0 9 5 9 4 5 23 24 4 8 24 7 6 1 14 8 3 8 14 4 23 12 8 7 22 5 7 8 1 2 4 10 2 6 3 8 19 3 4 9 20 4 6 12 14 3 4 1 12 3 -1 10 26 4 6 14 4 1 1 28 25 4 13 15 7 15 10 5 12 42 8 8 17 6 12 15 5 7 14 14 6 7 8 7 5 21 0 3 16 -2 14 18 22 16 2 13 5 22 17 5 17 6 10 6 12 0 13 17 6 7 5 12 9 2 39 4 15 5 7 -3 23 11 7 2 13 12 6 13 15 5 39 9 6 4 4 26 49 4 11 12 6 11 23 28 4 12 3 7 7 15 11 2 2 18 1 -3 2 9 17 5 8 15 10 34 7 39 13 6 0 8 7 5 16 7 20 7 6 3 12 49 16 13 -4 9 64 17 18 17 14 0 8 13 21 9 6 12 1 9 8 9 9 7 0 10 0 6 5 11 4 19 16 17 4 6 9 7 18 13 7 0 4 7 15 10 20 19 5 12 8 14 9 7 0 6 6 12 19 29 7 10 9 8 9 4 -3 14 2 19 6 9 4 16 6 20 6 17 8 8 8 2 7 30 -2 6 36 9 13 3 20 4 4 14 12 7 14 3 11 6 2 23 7 7 14 5 0 13 20 10 21 49 2 4 49 5 14 3 4 5 4 -1 13 9 15 10 8 17 15 9 18 35 5 7 13 7 19 12 4'
$wsMap['YDS'].Range('C2').Value = '4 5 28 0 9 2 4 11 2 2 -2 4 12 1 -3 1 3 1 1 3 3 3 3 3 4 -4 8 7 9 10 7 1 1 0 4 11 10 3 2 9 1 1 11 5 3 4 1 2 0 6 17 0 5 9 6 3 -5 6 3 -7 45 3 -3 4 4 -1 10 9 8 11 9 23 0 8 -5 3 0 3 3 -4 5 2 2 -3 0 6 1 11 3 12 4 5 -1 2 7 4 4 -3 7 0 4 3 2 13 6 4 8 6 2 -1 17 3 8 2 15 9 2 1 0 2 7 -3 1 -1 2 16 1 3 3 -4 4 4 7 -1 1 8 13 3 3 0 7 0 -2 2 7 3 2 3 3 2 -5 8 7 17 12 -4 5 6 2 9 2 3 2 8 14 5 3 9 12 4 4 3 -5 3 2 15 6 1 3 6 0 2 4 4 9 4 1 6 5 16 8 10 1 6 23 4 4 0 6 2 7 2 2 3 1 4 11 3 7 0 6 5 7 10 12 2 0 18 8 1 2 5 4 1 -5 2 0 0 5 0 1 1 4 -2 16 0 0 4 2 9 -1 3 43 3 7 9 7 16 6 0 37 5 6 1 2 0 4 3 10 11 6 2 5 0 7 0 3 0 0 9 2 3 2 4 9 3 4 13 13 -3 3 2 0 -3 3 5 3 2 10 0 -2 2 5 13 2 2 8 9 -1 1 1 2 1 3 1 9 4 1 0 1 1 17 2 2 6 2 10 -3 0 3 1 13 9 7 3 6 20 11 4 3 11 2 1 17 4 3 3 2 3 6 2 5 16 7 9 2 -3 0 4 9 8 3 5 6 3 8 13 7 4 6 6 3 4 11 4 9 23 1 4 14 18 7 3 3 1 6 2 8 -6 -1 5 -3 9 -1 6 6 8 2 7 3 2 10 3 8 10 4 2 8 4 6 0 3 8 7 2 5 12 -3 15 2 -1 3 4 0 4 6 -3 2 4 2 6 11 4 4 1 1 6 -1 7 4 6 5 3 6 2 8 6 6 -2 20 2 5 6 8 10 2 1 5 1 3 10 -1 1 -1 5 3 -3 2 7 3 5 4 9 3 7 8 0 6 1 1 1 2 4 3 3 6 5 2 8 3 7 5 2 2 3 8 6 0 16 3 6 2 3 9 5 3 9 3 10 -5 4 -2 2 4 9 3 9 2 2 1 -2 6 2 5 2 7 2 13 6 6 3 3 6 20 9 2 2 6 4 6 6 6 3 6 6 -2 1 11 8 2 6 -1 1 5 3 2 1 9 2 9 11 9 9 14 2 4 1 3 8 -1 3 1 6 -1 3 2 44 2 4 2 2 3 3 5 4 -3 4 2 -3 3 2 1 9 9 2 2 2 1 0 3 1 3 15 22 12 7 2 3 1 1 1 21 3 2 1 2 14 17 2 1 5 13 -1 7 -1 5 12 2 10 0 4 1 6 2 7 0 6 0 5 6 3 9 7 3 1 8 8 1 4 4 4 3 2 12 7 9 7 3 -1 10 2 10 4 4 3 4 2 2 6 6 5 4 3 1 3 4 9 -4 1 9 3 -1 11 8 -1 12 8 4 2 4 1 4 1 4 -5 -1 5 11 17 4 0 1 2 -2 2 4 2 5 5 3 2 2 17 2 3 6 0 3 10 0 5 6 12 8 13 13 1 3 5 1 1 16 3 3 2 6 2 2 2 7 1 3 2 3 14 5 -2 4 6 4 4 5 10 2 3 1 12 2 2 7 2 3 4 5 5 0 -3 4 4 7'
$wsMap['YDS'].Range('C3').Value = '15 9 9 6 7 3 13 19 3 5 6 5 1 3 7 9 3 10 7 16 20 3 10 10 17 4 37 4 18 15 8 10 10 9 11 9 -3 10 11 3 8 12 20 6 4 37 18 8 58 6 6 6 24 12 4 8 8 10 11 3 47 14 38 24 8 11 5 7 5 19 13 12 17 16 2 8 5 11 20 3 7 7 30 -8 4 6 3 4 37 6 2 28 27 8 6 9 29 19 2 17 8 14 15 19 9 3 2 13 8 8 23 21 6 6 5 -3 13 20 12 10 8 13 7 7 5 17 -1 29 14 1 12 57 2 8 12 23 22 3 1 17 14 3 -1 9 9 -2 4 4 36 19 18 14 9 35 19 8 6 11 11 17 1 2 35 9 2 19 12 49 3 16 23 1 13 36 1 9 10 18 20 16 8 29 14 4 18 15 6 7 -2 13 22 29 11 5 11 27 3 29 21 4 2 4 6 -7 5 -1 13 5 42 32 9 18 14 9 7 9 12 9 3 7 2 7 24 3 20 36 4 11 11 9 3 9 41 4 8 19 0 23 9 2 16 44 6 15 9 25 5 19 7 11 6 9 2 2 9 7 -1 16 7 11 2 3 1 13 11 15 19 11 12 11 15 -1 19 4 9 12 5 14 7 15 7 6 -3 21 0 18 1 8 7 3 7 18 11 18 10 5 10 6 17 9 2 4 6 9 6 6 5 6 10 16 2 3 4 5 7 9 8 4 1 5 4 12 5 8 6 11 13 -4 0 12 7 0 1 0 10 4 8 15 25 6 14 -8 14 8 2 7 10 7 10 20 14 32 5 26 16 24 4 17 8 5 11 9 13 4 18 13 46 13 15 6 10 31 5 4 13 9 9 5 11 2 9 4 11 12 13 26 0 4 31 2 9 17 5 19 16 9 25 12 21 9 12 11 29 5 14 23 20 25 12 8 46 7 8 12 14 18 47 4 30 3 6 2 10 11 2 18 13 28 8 4 19 3 9 -3 1 9 6 10 12 9 -1 5 8 25 7 19 9 0 23 20 25 20 16 16 9 5 8 1 9 5 26 17 21 3 8 1 2 9 10 26 5 12 6 14 8 38 19 8 18 12 9 7 4 16 0 13 13 11 10 12 3 7 5 16 4 9 10 4 1 6 9 3 19 -1 21 33 16 20 11 2 -2 11 23 10 17 24 19 4 30 11 14 1 13 3 6 14 4 5 -6 7 23 27 8 15 11 6 10 20 -1 9 17 5 2 10 10 0 7 2 15 9 7 19 5 8 11 5 14 6 6 9 18 27 9 11 1 7 10 17 16 4 4 17 -2 13 10 23 7 9 21 4 11 1 11 -1 9 2 19 5 14 7 19 11 20 4 15 -1 14 19 9 7 6 5 9 25 5 7 22 8 8 37 11 13 5 9 12 6 6 18 5 7 4 9 8 16 9 14 9 19 7 22 11 6 5 5 17 6 9 9 19 12 6 9 5 18 6 10 10 6 6 15 3 19 2 10 1 11 5 9 8 11 24 -1 21 -2 8 9 14 5 14 -5 21 7 11 26 20 3 8 15 5 11 0 6 3 15 7 4 3 25 36 3 0 15 6 17 12 27 1 9 1 7 11 5 7 18 1 5 11 14 6 6 6 36 5 3 7 11 15 22 10 23 9 8 14 6 18 7 9 3 19 19 13 18 11 6 25 6 28 3 36 9 9 30 3 15 13 21 15 3 4 8 1 13 0 4 8 -3 12 2 14 20 8 6 8 2 6 4 24 2 23 5 4 8 5 8'
$wsMap['ST'].Range('B4').Value = '61 60 65 70 60 67 72 66 65 67 57 64 65 65 66 69 63 66 67 67 68 66 60 66 66 63 58 61 67 62 65 65 60 63 65 64 55 70 67 65 64 61 61 66 67 57 66 52 60 65 66 64 59 60 64 53 69 45 52 67 66 64 63 63 66 64 66 64 64'
$wsMap['ST'].Range('B5').Value = '19 24 13 22 12 28 31 26 25 40 0 20 34 23 17 24 22 29 22 23 17 20 26 21 24 27 24 19 39 21 20 20 21 33 19 76 14 23 22 22 21 28 23 18 22 16 24 19 28 23 1 27 65 13 24 18 27 16 21 29 29 26 19 26 28 0 25 15 25'
$wsMap['ST'].Range('C4').Value = '17 20 25 14 12 13 26 24 23 19 29 20 19 16 23 18 22 22 27 26 25 4 13 19 14 45 27 24 10 16 24 32 22 29 16 24 21 25 19 21 25 28 21 16 17 27 18 24 22 26 19 30 38 20 26 33 24 24 23 32'
$wsMap['ST'].Range('D3').Value = '55 45 44 43 29 42 35 44 30 40 39 45 44 37 41 37 50 52 43 25 47 44 43 46 39 46 29 43 26 39 46 56 50 52 46 44 35 20 55 36 45 34 48 52 33 37 48 47 53 57 36 44 49 43 50 43 35 52 62 31 50 49 30 54 40 33 57 42 39 37 58 51 46 51 25 45 24 32 50 38 38 55 51 38 46 30 41 52 50 51 41 60 66 51 50 52 54 62 32 55 54 32 51 43 40 35 53 44 46 64 53 46 54 61'
$wsMap['ST'].Range('D4').Value = '15 0 0 0 0 3 0 0 0 0 0 6 0 0 11 0 -3 0 0 0 7 0 0 -1 0 0 0 8 0 0 2 0 4 14 0 6 0 0 14 0 9 0 8 1 0 5 0 0 21 6 0 0 0 11 9 0 0 2 11 0 7 -1 0 0 0 0 0 0 0 0 19 17 0 6 0 0 0 0 0 16 0 0 0 0 0 0 0 17 7 26 3 16 0 10 0 0 17 0 0 0 18 0 11 1 0 0 6 0 0 0 10 -2 19 0'
$wsMap['ST'].Range('D5').Value = '0 0 0 0 0 0 23 0 13 0'

# Update aggregate numeric totals
$wsMap['OFF'].Range('C2').Value = 388
$wsMap['OFF'].Range('F2').Value = 94
$wsMap['OFF'].Range('G2').Value = 88
$wsMap['OFF'].Range('J2').Value = 51
$wsMap['OFF'].Range('L2').Value = 572
$wsMap['OFF'].Range('M2').Value = 372
$wsMap['OFF'].Range('O2').Value = 46
$wsMap['OFF'].Range('P2').Value = 23
$wsMap['OFF'].Range('Q2').Value = 961
$wsMap['OFF'].Range('C3').Value = 313
$wsMap['OFF'].Range('E3').Value = 59
$wsMap['OFF'].Range('F3').Value = 241
$wsMap['OFF'].Range('G3').Value = 61
$wsMap['OFF'].Range('H3').Value = 57
$wsMap['OFF'].Range('I3').Value = 131
$wsMap['OFF'].Range('J3').Value = 111
$wsMap['OFF'].Range('N3').Value = 39
$wsMap['DEF'].Range('B2').Value = 15
$wsMap['DEF'].Range('C2').Value = 367
$wsMap['DEF'].Range('E2').Value = 18
$wsMap['DEF'].Range('F2').Value = 99
$wsMap['DEF'].Range('G2').Value = 116
$wsMap['DEF'].Range('I2').Value = 11
$wsMap['DEF'].Range('J2').Value = 56
$wsMap['DEF'].Range('L2').Value = 555
$wsMap['DEF'].Range('M2').Value = 384
$wsMap['DEF'].Range('O2').Value = 43
$wsMap['DEF'].Range('P2').Value = 20
$wsMap['DEF'].Range('Q2').Value = 962
$wsMap['DEF'].Range('B3').Value = 20
$wsMap['DEF'].Range('C3').Value = 375
$wsMap['DEF'].Range('E3').Value = 66
$wsMap['DEF'].Range('F3').Value = 215
$wsMap['DEF'].Range('G3').Value = 71
$wsMap['DEF'].Range('H3').Value = 46
$wsMap['DEF'].Range('I3').Value = 108
$wsMap['DEF'].Range('J3').Value = 111
$wsMap['ST'].Range('B2').Value = 149
$wsMap['ST'].Range('D2').Value = 114
$wsMap['ST'].Range('F2').Value = 88
$wsMap['ST'].Range('G2').Value = 84
$wsMap['ST'].Range('L2').Value = 28
$wsMap['ST'].Range('M2').Value = 21
$wsMap['ST'].Range('N2').Value = 13
$wsMap['ST'].Range('O2').Value = 12
$wsMap['ST'].Range('B3').Value = 80
$wsMap['TURNS'].Range('C2').Value = 8
$wsMap['TURNS'].Range('D2').Value = 19
$wsMap['TURNS'].Range('D3').Value = 18
$wsMap['TURNS'].Range('E3').Value = 13
$wsMap['PEN'].Range('D2').Value = 15
$wsMap['PEN'].Range('B3').Value = 21
